$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.186.63"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "1.601.27"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "212.03"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "0.481"
$ws.Range("E7").Value = "  +1.08%  "
$ws.Range("D8").Value = "0.249"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "18.23"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").Value = "0.0809"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").Value = "1.823.17"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "1.596.42"
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("D15").Value = "0.518"
$ws.Range("E15").Value = "  +2.32%  "
$ws.Range("D16").Value = "26.183.39"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "60.96"
$ws.Range("E17").Value = "  +1.23%  "
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "204.06"
$ws.Range("E20").Value = "  +4.39%  "
$ws.Range("D21").Value = "4.27"
$ws.Range("E21").Value = "  +1.44%  "
$ws.Range("D22").Value = "9.29"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("E24").Value = "  +12.58%  "
$ws.Range("D25").Value = "143.36"
$ws.Range("E25").Value = "  +1.08%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -7.67%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  +1.12%  "
$ws.Range("D31").Value = "1.17"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").Value = "1.129.41"
$ws.Range("E36").Value = "  +2.36%  "
$ws.Range("E37").Value = "  +7.56%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "0.793"
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("E41").Value = "  -1.92%  "
$ws.Range("D42").Value = "0.781"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "5.17"
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "1.736.94"
$ws.Range("E44").Value = "  +0.73%  "
$ws.Range("D45").Value = "92.02"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("D47").Value = "54.19"
$ws.Range("E47").Value = "  +1.53%  "
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("D49").Value = "0.406"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("B50").Value = "USDD"
$ws.Range("C50").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0952"
$ws.Range("E51").Value = "  -13.89%  "
